$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 15:34"

# Row 119/120: Islandia overtakes Zambia in the ranking (table sorted desc by
# total cases), so the two countries swap rows; row 119 gets Islandia's
# refreshed numbers and row 120 gets Zambia's (unchanged) numbers.
$ws.Range("A119").Value = "Islandia"
$ws.Range("A120").Value = "Zambia"

# Country statistics refresh
$ws.Range("B4").Value = 3357130
$ws.Range("C4").Value = 1484
$ws.Range("E4").Value = 1729010
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 137418
$ws.Range("B6").Value = 856062
$ws.Range("C6").Value = 5704
$ws.Range("D6").Value = 541227
$ws.Range("E6").Value = 292073
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = 22762
$ws.Range("B17").Value = 232259
$ws.Range("C17").Value = 2779
$ws.Range("D17").Value = 167138
$ws.Range("E17").Value = 62898
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 2223
$ws.Range("B19").Value = 199828
$ws.Range("C19").Value = 16
$ws.Range("E19").Value = 6194
$ws.Range("B24").Value = 103598
$ws.Range("C24").Value = 470
$ws.Range("D24").Value = 99743
$ws.Range("E24").Value = 3708
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 147
$ws.Range("B40").Value = 51022
$ws.Range("C40").Value = 101
$ws.Range("B63").Value = 18360
$ws.Range("C63").Value = 287
$ws.Range("D63").Value = 13876
$ws.Range("E63").Value = 4091
$ws.Range("G63").Value = 11
$ws.Range("H63").Value = 393
$ws.Range("B76").Value = 9674
$ws.Range("C76").Value = 283
$ws.Range("D76").Value = 5634
$ws.Range("E76").Value = 3780
$ws.Range("B78").Value = 8979
$ws.Range("C78").Value = 2
$ws.Range("E78").Value = 589
$ws.Range("B84").Value = 7294
$ws.Range("C84").Value = 3
$ws.Range("E84").Value = 165
$ws.Range("B119").Value = 1896
$ws.Range("C119").Value = 8
$ws.Range("D119").Value = 1865
$ws.Range("E119").Value = 21
$ws.Range("H119").Value = 10
$ws.Range("B120").Value = 1895
$ws.Range("D120").Value = 1348
$ws.Range("E120").Value = 505
$ws.Range("H120").Value = 42
$ws.Range("D123").Value = 1469
$ws.Range("E123").Value = 261
$ws.Range("B142").Value = 1010
$ws.Range("C142").Value = 12
$ws.Range("D142").Value = 423
$ws.Range("E142").Value = 536
$ws.Range("G142").Value = 4
$ws.Range("H142").Value = 51
$ws.Range("B158").Value = 483
$ws.Range("C158").Value = 21
$ws.Range("E158").Value = 340
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 25
